$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

$wsZhCn.Range("H7").Value = "2016-08-17 12:39:17"
$wsDeDe.Range("H7").Value = "2016-08-17 12:39:23"
$wsOverview.Range("G7").Value = "2016-08-17 12:39:23"
